# e2e-document.xlsx: fill in the "Result" column for the first test sheet
# ("City search") with the actual e2e test results, mark those cells with
# the built-in "Good" cell style, and leave the first sheet as the active
# tab/selection (second sheet just gets its selection moved off A6).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Row 2: "no city entered + search" -> returns every record in the DB
$ws1.Range("C2").Value = "מחזיר את כלל התוצאות במאגר כנדרש"

# Row 3: "city entered + search" -> entered 'Holon', got a single matching result
$ws1.Range("C3").Value = "הוזן 'חולון', התקלה תוצאה 1 כנדרש"

# Row 4: "different city entered + search" -> entered 'Tel Aviv', Holon result
# is gone, only Tel Aviv results come back
$ws1.Range("C4").Value = "הוזן 'תל אביב', התוצאה `nשל חולון לא חזרה והתקבלו תוצאות מתל אביב בלבד כנדרש"

# Mark the newly-filled results with the built-in "Good" (green) cell style
$ws1.Range("C2:C4").Style = "Good"
$ws1.Range("C2:C4").WrapText = $true

# Move sheet2's selection off A6, then land the final selection on sheet1,
# which makes sheet1 the active tab/sheet saved with the workbook.
$ws2.Range("A2:C2").Select()
$ws1.Range("C5").Select()
